$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
Write-Host "Slide 6 shapes count: $($s.Shapes.Count)"
